# Update "想去人数" (F column) figures on both the "展览" and "全部类型"
# sheets to match the freshly scraped counts.
$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 2986
    5  = 6733
    6  = 1729
    10 = 121
    11 = 6
    13 = 132
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
